$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column (BT = col 72) uses the same column width as its neighbours (raw
# OOXML width="12", which corresponds to a COM ColumnWidth of 11.17 given this
# sheet's base font).
$ws.Columns.Item(72).ColumnWidth = 11.17

# Reference cells carrying the three distinct cell styles already used on this
# sheet (no fill / yellow fill / light-blue fill), used as Copy() sources so the
# new BT column cells land on the SAME style entries instead of minting new ones.
$styleRefs = @{ "1" = $ws.Range("A2"); "2" = $ws.Range("D2"); "3" = $ws.Range("N2") }

# --- Row 1 header: date label (literal text, not an Excel date) ---
$cell = $ws.Cells.Item(1, 72)
$styleRefs["1"].Copy($cell)
$cell.NumberFormat = "@"
$cell.Value = "2024/11/19"

# --- Rows 2-53: numeric data for 2024/11/19 ---
$data = @(
    @{ Row = 2; Style = "1"; Value = 205.2 },
    @{ Row = 3; Style = "3"; Value = 137.1 },
    @{ Row = 4; Style = "1"; Value = 146 },
    @{ Row = 5; Style = "1"; Value = 144.3 },
    @{ Row = 6; Style = "3"; Value = 131.9 },
    @{ Row = 7; Style = "1"; Value = 190.6 },
    @{ Row = 8; Style = "1"; Value = 178.6 },
    @{ Row = 9; Style = "1"; Value = 332.6 },
    @{ Row = 10; Style = "1"; Value = 261.7 },
    @{ Row = 11; Style = "1"; Value = 165.9 },
    @{ Row = 12; Style = "1"; Value = 146.3 },
    @{ Row = 13; Style = "1"; Value = 237.3 },
    @{ Row = 14; Style = "1"; Value = 161.2 },
    @{ Row = 15; Style = "1"; Value = 202.8 },
    @{ Row = 16; Style = "3"; Value = 136.2 },
    @{ Row = 17; Style = "1"; Value = 167.9 },
    @{ Row = 18; Style = "1"; Value = 178.4 },
    @{ Row = 19; Style = "2"; Value = 105.4 },
    @{ Row = 20; Style = "1"; Value = 191.9 },
    @{ Row = 21; Style = "1"; Value = 144.2 },
    @{ Row = 22; Style = "1"; Value = 197.1 },
    @{ Row = 23; Style = "1"; Value = 157.5 },
    @{ Row = 24; Style = "1"; Value = 170.3 },
    @{ Row = 25; Style = "3"; Value = 125 },
    @{ Row = 26; Style = "3"; Value = 134.7 },
    @{ Row = 27; Style = "1"; Value = 182.2 },
    @{ Row = 28; Style = "1"; Value = 153.6 },
    @{ Row = 29; Style = "1"; Value = 169.4 },
    @{ Row = 30; Style = "1"; Value = 202.2 },
    @{ Row = 31; Style = "1"; Value = 166.9 },
    @{ Row = 32; Style = "2"; Value = 103.8 },
    @{ Row = 33; Style = "1"; Value = 166.3 },
    @{ Row = 34; Style = "1"; Value = 166.5 },
    @{ Row = 35; Style = "1"; Value = 190 },
    @{ Row = 36; Style = "1"; Value = 184.6 },
    @{ Row = 37; Style = "1"; Value = 159.2 },
    @{ Row = 38; Style = "1"; Value = 152.2 },
    @{ Row = 39; Style = "1"; Value = 163.1 },
    @{ Row = 40; Style = "1"; Value = 184.4 },
    @{ Row = 41; Style = "1"; Value = 155.4 },
    @{ Row = 42; Style = "1"; Value = 228.6 },
    @{ Row = 43; Style = "1"; Value = 217.2 },
    @{ Row = 44; Style = "1"; Value = 146.8 },
    @{ Row = 45; Style = "2"; Value = 123.9 },
    @{ Row = 46; Style = "1"; Value = 197.1 },
    @{ Row = 47; Style = "3"; Value = 138.7 },
    @{ Row = 48; Style = "1"; Value = 170.6 },
    @{ Row = 49; Style = "1"; Value = 264.3 },
    @{ Row = 50; Style = "2"; Value = 117 },
    @{ Row = 51; Style = "1"; Value = 165.6 },
    @{ Row = 52; Style = "3"; Value = 133.8 },
    @{ Row = 53; Style = "1"; Value = 144.4 }
)

foreach ($item in $data) {
    $cell = $ws.Cells.Item($item.Row, 72)
    $styleRefs[$item.Style].Copy($cell)
    $cell.Value = $item.Value
}

Write-Output "BT column (2024/11/19) written"